# ADD results from server
# Updates the single data row (row 2) of the "2025", "2030", "2035",
# "2040" and "2045" sheets with refreshed values received from the server.
# The "2050" sheet is unchanged.

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 0
$ws2025.Range("B2").Value = 0.003837439598569248
$ws2025.Range("E2").Value = 0.3718167378372543
$ws2025.Range("G2").Value = 0.2494892361374791
$ws2025.Range("I2").Value = 0.368798651279322
$ws2025.Range("L2").Value = 0.597153
$ws2025.Range("M2").Value = 0.0822565
$ws2025.Range("N2").Value = 12.82009457445623
$ws2025.Range("O2").Value = 3.537984783585709

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 0.009688525212652177
$ws2030.Range("B2").Value = 0.04891593573387359
$ws2030.Range("E2").Value = 0.221684650692804
$ws2030.Range("I2").Value = 0.4222272367131874
$ws2030.Range("L2").Value = 0.1141686999999999
$ws2030.Range("M2").Value = 0.04737166666666669
$ws2030.Range("N2").Value = 4.970917217357455
$ws2030.Range("O2").Value = 2.349741755853875

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 0.03875691745240942
$ws2035.Range("B2").Value = 0.03849685245042471
$ws2035.Range("E2").Value = 0.198894071362835
$ws2035.Range("I2").Value = 0.4598544380751153
$ws2035.Range("L2").Value = 0
$ws2035.Range("M2").Value = 0.05372516997467353
$ws2035.Range("N2").Value = 9.039033233069347
$ws2035.Range("O2").Value = 4.931912879036898

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 0.001072054624113319

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 0.159740130277202
